$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix naive component forecaster bug: clear stale early-period forecast values
# that should not have been populated (insufficient history).
$cellsToClear = @("C2", "E2", "C3", "C4")
foreach ($cellRef in $cellsToClear) {
    $ws.Range($cellRef).ClearContents()
}

# Updated forecast values resulting from the corrected naive forecaster computation
# (small floating-point differences from the recalculated averages).
$cellUpdates = @{
    "E3" = 9.591339540850829
    "E4" = 4.422525088127305
    "C6" = -14.45332333832744
    "E6" = -2.928447329610051
    "E7" = -2.225127715916664
    "C8" = 8.600536527919612
    "C9" = 9.399485634179205
    "C11" = 5.169490031659651
    "E11" = 9.213376886330327
    "C12" = 4.639893381363192
    "C13" = -0.3722371047999995
    "E13" = 2.684220738731979
    "E14" = 2.429116709932599
    "C15" = 4.098801479368319
    "E16" = 3.941300050092877
    "E17" = 2.714258593289998
    "C19" = 2.352205130086094
    "C21" = 4.083548352538391
    "E21" = 3.58625614607444
    "C22" = 4.695933104194361
    "C24" = 4.861590900330715
    "C25" = 5.402237127943765
    "E25" = 4.104053120889195
    "C29" = 0.8513583007189407
    "E29" = 2.225279621195853
    "C31" = 1.015697339178057
    "E31" = 2.122104735451624
    "E32" = -0.6322362079330235
    "E34" = -1.352810423674367
    "C35" = 4.074459326939817
    "E35" = -0.2414327668618488
    "E36" = 1.077755602068309
    "E37" = 1.148476797857967
    "E38" = -0.3934198590721305
    "E41" = 1.534407168230811
    "C42" = 5.120680133083622
    "C43" = 5.356482122456163
    "E43" = 12.6296844023545
    "C46" = -0.5532735011319123
    "C47" = -2.464475897442031
    "E48" = 1.793234865396331
    "C49" = -0.8995735674421024
    "E49" = 0.3402056885013494
    "E50" = -1.28528149926006
    "E51" = -1.655020334777801
    "E52" = -0.6714033493142035
    "E53" = -0.5945514555738662
}
foreach ($cellRef in $cellUpdates.Keys) {
    $ws.Range($cellRef).Value = $cellUpdates[$cellRef]
}
